{"js": "const replacements = [\n  [\"97\u00d774=\", \"54\u00d748=\"],\n  [\"96\u00d770=\", \"22\u00d774=\"],\n  [\"76\u00d736=\", \"44\u00d794=\"],\n  [\"43\u00d751=\", \"96\u00d784=\"],\n  [\"36\u00d785=\", \"46\u00d729=\"],\n  [\"59\u00d741=\", \"19\u00d735=\"],\n  [\"81\u00d742=\", \"41\u00d712=\"],\n  [\"32\u00d770=\", \"66\u00d761=\"],\n  [\"28\u00d761=\", \"31\u00d756=\"],\n  [\"32\u00d771=\", \"81\u00d771=\"],\n  [\"32\u00d796=\", \"36\u00d793=\"],\n  [\"84\u00d787=\", \"57\u00d788=\"],\n  [\"92\u00d784=\", \"37\u00d731=\"],\n  [\"81\u00d780=\", \"23\u00d791=\"],\n  [\"71\u00d780=\", \"11\u00d797=\"],\n  [\"86\u00d731=\", \"33\u00d797=\"],\n  [\"28\u00d780=\", \"68\u00d753=\"],\n  [\"99\u00d753=\", \"40\u00d767=\"],\n  [\"33\u00d732=\", \"97\u00d781=\"],\n  [\"61\u00d793=\", \"65\u00d754=\"],\n  [\"86\u00d732=\", \"41\u00d730=\"],\n  [\"31\u00d791=\", \"66\u00d754=\"],\n  [\"96\u00d791=\", \"21\u00d713=\"],\n  [\"75\u00d715=\", \"82\u00d754=\"],\n  [\"31\u00d788=\", \"77\u00d797=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"97\u00d774=\", \"54\u00d748=\"),\n    @(\"96\u00d770=\", \"22\u00d774=\"),\n    @(\"76\u00d736=\", \"44\u00d794=\"),\n    @(\"43\u00d751=\", \"96\u00d784=\"),\n    @(\"36\u00d785=\", \"46\u00d729=\"),\n    @(\"59\u00d741=\", \"19\u00d735=\"),\n    @(\"81\u00d742=\", \"41\u00d712=\"),\n    @(\"32\u00d770=\", \"66\u00d761=\"),\n    @(\"28\u00d761=\", \"31\u00d756=\"),\n    @(\"32\u00d771=\", \"81\u00d771=\"),\n    @(\"32\u00d796=\", \"36\u00d793=\"),\n    @(\"84\u00d787=\", \"57\u00d788=\"),\n    @(\"92\u00d784=\", \"37\u00d731=\"),\n    @(\"81\u00d780=\", \"23\u00d791=\"),\n    @(\"71\u00d780=\", \"11\u00d797=\"),\n    @(\"86\u00d731=\", \"33\u00d797=\"),\n    @(\"28\u00d780=\", \"68\u00d753=\"),\n    @(\"99\u00d753=\", \"40\u00d767=\"),\n    @(\"33\u00d732=\", \"97\u00d781=\"),\n    @(\"61\u00d793=\", \"65\u00d754=\"),\n    @(\"86\u00d732=\", \"41\u00d730=\"),\n    @(\"31\u00d791=\", \"66\u00d754=\"),\n    @(\"96\u00d791=\", \"21\u00d713=\"),\n    @(\"75\u00d715=\", \"82\u00d754=\"),\n    @(\"31\u00d788=\", \"77\u00d797=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
